# chore: update Sheets via scheduled runner
# Refresh cached market-board price/profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit* columns) across all profession leve-profit tables.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 548.1111
$ws.Cells.Item(11, 9).Value = 548.1111
$ws.Cells.Item(11, 11).Value = 548.1111
$ws.Cells.Item(11, 13).Value = -408.1111
$ws.Cells.Item(40, 8).Value = 1920
$ws.Cells.Item(86, 8).Value = 3480.6667
$ws.Cells.Item(86, 9).Value = 2757.625
$ws.Cells.Item(86, 10).Value = 4059.1
$ws.Cells.Item(86, 11).Value = 2757.625
$ws.Cells.Item(86, 12).Value = 4059.1
$ws.Cells.Item(86, 13).Value = -1634.625
$ws.Cells.Item(86, 14).Value = -6305.1
$ws.Cells.Item(89, 8).Value = 3480.6667
$ws.Cells.Item(89, 9).Value = 2757.625
$ws.Cells.Item(89, 10).Value = 4059.1
$ws.Cells.Item(89, 11).Value = 13788.125
$ws.Cells.Item(89, 12).Value = 20295.5
$ws.Cells.Item(89, 13).Value = -8172.125
$ws.Cells.Item(89, 14).Value = -31527.5
$ws.Cells.Item(106, 8).Value = 32992.43
$ws.Cells.Item(106, 9).Value = 32992.43
$ws.Cells.Item(106, 11).Value = 32992.43
$ws.Cells.Item(106, 13).Value = -32361.43
$ws.Cells.Item(116, 8).Value = 9000
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 13).ClearContents()
$ws.Cells.Item(127, 8).Value = 1433.3334
$ws.Cells.Item(127, 9).Value = 1150
$ws.Cells.Item(127, 11).Value = 3450
$ws.Cells.Item(127, 13).Value = 1510
$ws.Cells.Item(138, 8).Value = 11728.85
$ws.Cells.Item(138, 10).Value = 12461.311
$ws.Cells.Item(138, 12).Value = 37383.933
$ws.Cells.Item(138, 14).Value = -47663.933

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 20958.115
$ws.Cells.Item(32, 9).Value = 15264.286
$ws.Cells.Item(32, 10).Value = 29498.857
$ws.Cells.Item(32, 11).Value = 15264.286
$ws.Cells.Item(32, 12).Value = 29498.857
$ws.Cells.Item(32, 13).Value = -14977.286
$ws.Cells.Item(32, 14).Value = -30072.857
$ws.Cells.Item(102, 8).Value = 2133
$ws.Cells.Item(102, 9).Value = 2190.2
$ws.Cells.Item(102, 10).Value = 2101.2222
$ws.Cells.Item(102, 11).Value = 2190.2
$ws.Cells.Item(102, 12).Value = 2101.2222
$ws.Cells.Item(102, 13).Value = -568.1999999999998
$ws.Cells.Item(102, 14).Value = -5345.2222

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5, 8).Value = 2212.875
$ws.Cells.Item(5, 9).Value = 1700.5
$ws.Cells.Item(5, 10).Value = 3750
$ws.Cells.Item(5, 11).Value = 1700.5
$ws.Cells.Item(5, 12).Value = 3750
$ws.Cells.Item(5, 13).Value = -1587.5
$ws.Cells.Item(5, 14).Value = -3976
$ws.Cells.Item(105, 8).Value = 4127.4443
$ws.Cells.Item(105, 9).Value = 3374.5557
$ws.Cells.Item(105, 11).Value = 3374.5557
$ws.Cells.Item(105, 13).Value = -1627.5557
$ws.Cells.Item(106, 8).Value = 27967.5
$ws.Cells.Item(106, 10).Value = 27967.5
$ws.Cells.Item(106, 12).Value = 27967.5
$ws.Cells.Item(106, 14).Value = -30491.5
$ws.Cells.Item(134, 8).Value = 2470.2173
$ws.Cells.Item(134, 9).Value = 1472.2778
$ws.Cells.Item(134, 11).Value = 4416.8334
$ws.Cells.Item(134, 13).Value = -1881.8334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 245.75
$ws.Cells.Item(22, 9).Value = 138
$ws.Cells.Item(22, 11).Value = 138
$ws.Cells.Item(22, 13).Value = 212
$ws.Cells.Item(31, 8).Value = 6505.923
$ws.Cells.Item(31, 9).Value = 5211.25
$ws.Cells.Item(31, 11).Value = 5211.25
$ws.Cells.Item(31, 13).Value = -4916.25
$ws.Cells.Item(34, 8).Value = 6505.923
$ws.Cells.Item(34, 9).Value = 5211.25
$ws.Cells.Item(34, 11).Value = 5211.25
$ws.Cells.Item(34, 13).Value = -5009.25
$ws.Cells.Item(68, 8).Value = 44499.5
$ws.Cells.Item(68, 10).Value = 44499.5
$ws.Cells.Item(68, 12).Value = 44499.5
$ws.Cells.Item(68, 14).Value = -45997.5
$ws.Cells.Item(71, 8).Value = 44499.5
$ws.Cells.Item(71, 10).Value = 44499.5
$ws.Cells.Item(71, 12).Value = 133498.5
$ws.Cells.Item(71, 14).Value = -140986.5
$ws.Cells.Item(105, 8).Value = 5701.5557
$ws.Cells.Item(105, 9).Value = 4495.75
$ws.Cells.Item(105, 10).Value = 6666.2
$ws.Cells.Item(105, 11).Value = 4495.75
$ws.Cells.Item(105, 12).Value = 6666.2
$ws.Cells.Item(105, 13).Value = -2748.75
$ws.Cells.Item(105, 14).Value = -10160.2
$ws.Cells.Item(107, 8).Value = 925.0714
$ws.Cells.Item(107, 9).Value = 580.5714
$ws.Cells.Item(107, 10).Value = 1269.5714
$ws.Cells.Item(107, 11).Value = 580.5714
$ws.Cells.Item(107, 12).Value = 1269.5714
$ws.Cells.Item(107, 13).Value = 1339.4286
$ws.Cells.Item(107, 14).Value = -5109.5714
$ws.Cells.Item(132, 8).Value = 2344.4285
$ws.Cells.Item(132, 9).Value = 2303.963
$ws.Cells.Item(132, 11).Value = 6911.889000000001
$ws.Cells.Item(132, 13).Value = -4381.889000000001
$ws.Cells.Item(141, 8).Value = 49833.332
$ws.Cells.Item(141, 10).Value = 49833.332
$ws.Cells.Item(141, 12).Value = 49833.332
$ws.Cells.Item(141, 14).Value = -60193.332

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 433.45
$ws.Cells.Item(5, 10).Value = 1000
$ws.Cells.Item(5, 12).Value = 3000
$ws.Cells.Item(5, 14).Value = -3224
$ws.Cells.Item(32, 8).Value = 3885830.2
$ws.Cells.Item(32, 10).Value = 4113879
$ws.Cells.Item(32, 12).Value = 12341637
$ws.Cells.Item(32, 14).Value = -12342203
$ws.Cells.Item(34, 8).Value = 3833.4443
$ws.Cells.Item(34, 9).Value = 2425.5
$ws.Cells.Item(34, 10).Value = 4959.8
$ws.Cells.Item(34, 11).Value = 7276.5
$ws.Cells.Item(34, 12).Value = 14879.4
$ws.Cells.Item(34, 13).Value = -7192.5
$ws.Cells.Item(34, 14).Value = -15047.4
$ws.Cells.Item(39, 8).Value = 2250
$ws.Cells.Item(39, 9).Value = 1666.6666
$ws.Cells.Item(39, 11).Value = 4999.9998
$ws.Cells.Item(39, 13).Value = -4705.9998
$ws.Cells.Item(55, 8).Value = 101649.7
$ws.Cells.Item(55, 9).Value = 250749.25
$ws.Cells.Item(55, 10).Value = 2250
$ws.Cells.Item(55, 11).Value = 752247.75
$ws.Cells.Item(55, 12).Value = 6750
$ws.Cells.Item(55, 13).Value = -752070.75
$ws.Cells.Item(55, 14).Value = -7104
$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(129, 14).ClearContents()
$ws.Cells.Item(135, 8).Value = 433.45
$ws.Cells.Item(135, 10).Value = 1000
$ws.Cells.Item(135, 12).Value = 9000
$ws.Cells.Item(135, 14).Value = -14070

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 2146.5557
$ws.Cells.Item(97, 9).Value = 2246.2856
$ws.Cells.Item(97, 11).Value = 2246.2856
$ws.Cells.Item(97, 13).Value = -1750.2856
$ws.Cells.Item(122, 8).Value = 397778.44
$ws.Cells.Item(122, 10).Value = 718439.2
$ws.Cells.Item(122, 12).Value = 2155317.6
$ws.Cells.Item(122, 14).Value = -2160217.6
$ws.Cells.Item(126, 8).Value = 4082.4614
$ws.Cells.Item(126, 9).Value = 3509
$ws.Cells.Item(126, 11).Value = 10527
$ws.Cells.Item(126, 13).Value = -8057
$ws.Cells.Item(132, 8).Value = 4221.769
$ws.Cells.Item(132, 9).Value = 3226.8
$ws.Cells.Item(132, 10).Value = 4843.625
$ws.Cells.Item(132, 11).Value = 9680.400000000001
$ws.Cells.Item(132, 12).Value = 14530.875
$ws.Cells.Item(132, 13).Value = -7150.400000000001
$ws.Cells.Item(132, 14).Value = -19590.875
$ws.Cells.Item(141, 8).Value = 57357.25
$ws.Cells.Item(141, 10).Value = 57357.25
$ws.Cells.Item(141, 12).Value = 57357.25
$ws.Cells.Item(141, 14).Value = -67717.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 2989.6
$ws.Cells.Item(55, 9).Value = 2737
$ws.Cells.Item(55, 10).Value = 4000
$ws.Cells.Item(55, 11).Value = 2737
$ws.Cells.Item(55, 12).Value = 4000
$ws.Cells.Item(55, 13).Value = -2564
$ws.Cells.Item(55, 14).Value = -4346
$ws.Cells.Item(68, 8).Value = 2600.25
$ws.Cells.Item(68, 9).Value = 2000
$ws.Cells.Item(68, 10).Value = 2800.3333
$ws.Cells.Item(68, 11).Value = 2000
$ws.Cells.Item(68, 12).Value = 2800.3333
$ws.Cells.Item(68, 13).Value = -1251
$ws.Cells.Item(68, 14).Value = -4298.3333
$ws.Cells.Item(71, 8).Value = 2600.25
$ws.Cells.Item(71, 9).Value = 2000
$ws.Cells.Item(71, 10).Value = 2800.3333
$ws.Cells.Item(71, 11).Value = 10000
$ws.Cells.Item(71, 12).Value = 14001.6665
$ws.Cells.Item(71, 13).Value = -6256
$ws.Cells.Item(71, 14).Value = -21489.6665
$ws.Cells.Item(82, 8).Value = 2779.1
$ws.Cells.Item(82, 9).Value = 2723.875
$ws.Cells.Item(82, 10).Value = 3000
$ws.Cells.Item(82, 11).Value = 2723.875
$ws.Cells.Item(82, 12).Value = 3000
$ws.Cells.Item(82, 13).Value = -2362.875
$ws.Cells.Item(82, 14).Value = -3722
$ws.Cells.Item(85, 8).Value = 2779.1
$ws.Cells.Item(85, 9).Value = 2723.875
$ws.Cells.Item(85, 10).Value = 3000
$ws.Cells.Item(85, 11).Value = 2723.875
$ws.Cells.Item(85, 12).Value = 3000
$ws.Cells.Item(85, 13).Value = -1475.875
$ws.Cells.Item(85, 14).Value = -5496
$ws.Cells.Item(132, 8).Value = 4280.08
$ws.Cells.Item(132, 9).Value = 3725.1
$ws.Cells.Item(132, 11).Value = 11175.3
$ws.Cells.Item(132, 13).Value = -8645.299999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3992.2856
$ws.Cells.Item(81, 9).Value = 3684.077
$ws.Cells.Item(81, 10).Value = 7999
$ws.Cells.Item(81, 11).Value = 7368.154
$ws.Cells.Item(81, 12).Value = 15998
$ws.Cells.Item(81, 13).Value = -6307.154
$ws.Cells.Item(81, 14).Value = -18120
$ws.Cells.Item(84, 8).Value = 3992.2856
$ws.Cells.Item(84, 9).Value = 3684.077
$ws.Cells.Item(84, 10).Value = 7999
$ws.Cells.Item(84, 11).Value = 36840.77
$ws.Cells.Item(84, 12).Value = 79990
$ws.Cells.Item(84, 13).Value = -31536.77
$ws.Cells.Item(84, 14).Value = -90598
$ws.Cells.Item(101, 8).Value = 25000
$ws.Cells.Item(101, 10).Value = 25000
$ws.Cells.Item(101, 12).Value = 25000
$ws.Cells.Item(101, 14).Value = -31490
$ws.Cells.Item(122, 8).Value = 2332.6667
$ws.Cells.Item(122, 9).Value = 2399.2
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 7197.599999999999
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -4747.599999999999
$ws.Cells.Item(122, 14).Value = -10900
$ws.Cells.Item(126, 8).Value = 126786.375
$ws.Cells.Item(126, 10).Value = 2100
$ws.Cells.Item(126, 12).Value = 6300
$ws.Cells.Item(126, 14).Value = -11240
$ws.Cells.Item(136, 8).Value = 78929.62
$ws.Cells.Item(136, 9).Value = 1553.7273
$ws.Cells.Item(136, 11).Value = 4661.1819
$ws.Cells.Item(136, 13).Value = -2111.1819
